$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H5").Value = "2016-08-26 20:40:49"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H5").Value = "2016-08-26 20:40:53"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G5").Value = "2016-08-26 20:40:53"
